$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Closing / opening balance columns (G, H) for the data rows now hold
# ratio-style text ("2000 : 1" / "500 : 1") instead of raw numbers.
$ws.Range("G2:G4").Value = "2000 : 1"
$ws.Range("H2:H4").Value = "500 : 1"

# Tidy up the now-unused placeholder cells in column A (rows 5-7 have no
# content or formatting, Excel drops the empty <c> nodes on save).
$ws.Range("A5:A7").ClearContents()

# Reflect the new selection left behind after editing the balance cells.
$ws.Range("G2:H4").Select()
